$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cases query (B2): append an ORDER BY / LIMIT clause -------------------
$casesQuery = $ws.Range("B2").Value2
$casesQuery = $casesQuery + "`n order By ss.study_subject_id ASC LIMIT 100 "
$ws.Range("B2").Value2 = $casesQuery

# --- Samples query (B3): append an ORDER BY / LIMIT clause -----------------
$samplesQuery = $ws.Range("B3").Value2
$samplesQuery = $samplesQuery + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Value2 = $samplesQuery

# --- Files query (B4): replace the lowercase "order by" clause -------------
$filesQuery = $ws.Range("B4").Value2
$filesQuery = $filesQuery.Replace("    order by f.file_name", "    order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value2 = $filesQuery

# The extra line in each query makes the wrapped cell text taller, so the
# (wrap-text) rows auto-grow to fit the new content.
$ws.Rows(2).RowHeight = 331.2
$ws.Rows(3).RowHeight = 360

# Update the on-screen selection/scroll position left by the last save.
$ws.Range("C3").Select()
